# Commit: "Added Display class Updated images"
#
# The TimeRelay sheet previously had two boolean-ish columns, "IsTopBuyed"
# (col G) and "IsNew" (col H), neither of which held any data in the body
# rows. This edit removes the "IsNew" column entirely and repurposes the
# "IsTopBuyed" column header into a new "DisplayClass" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TimeRelay")

# Delete column H ("IsNew") entirely - this shifts every column after it
# one place to the left (old I..V become H..U).
$ws.Range("H1").EntireColumn.Delete()

# Rename what used to be column G's header ("IsTopBuyed") to the new
# "DisplayClass" header.
$ws.Range("G1").Value = "DisplayClass"

# Update the active selection on the TimeRelay sheet to G2 (previously G10).
$ws.Range("G2").Select()
